$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "ID" values for a handful of workers in column C.
# New unique values are entered in this order so the shared-string table
# is populated identically to the source edit (C26, C35, C14, C11, C7).
$ws.Range("C26").Value = "F72F821E"
$ws.Range("C35").Value = "FA02F7A3"
$ws.Range("C14").Value = "50483E8D"
$ws.Range("C11").Value = "99366BB0"
$ws.Range("C7").Value = "F90B3594"

# C7, C14 and C35 previously held the "empty ID, bordered" style; once they
# get a real ID the formatting should match the other filled-in ID cells
# (left-aligned custom font, same as column A's alternating style). Copy
# that formatting from the matching A-column cell on each row.
$ws.Range("A7").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Range("A14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("A35").Copy()
$ws.Range("C35").PasteSpecial(-4122)

# Move the active selection to C2, matching where the user left off editing.
$ws.Range("C2").Select()
